# [PHOENIX-6081] completed official register complaint with flow
# Adds a new "processingStatus" / "PROCESSING" column (G) to the
# grievanceDetails sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell G1 = "processingStatus"
$ws.Range("G1").Value = "processingStatus"
# New data cell G2 = "PROCESSING"
$ws.Range("G2").Value = "PROCESSING"

# Header G1 should use the same bold/Courier-New/vertical-centered format
# already used by the other header-style cell (A2), so copy that cell's
# formatting (not its value) onto G1.
$ws.Range("A2").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen the new column like the rest of the sheet's columns.
$ws.Columns.Item(7).ColumnWidth = 19.140625

# Match the new selection recorded in the sheet view.
$ws.Range("L7").Select()
